# Gunshot Detection Project Poster - fill in 1D CNN Keras/TFLite results
# and swap the mislabeled Precision/Accuracy rows in the results table.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Table 32" is the results table with columns:
#   1:(row labels) 2:1D-Keras 3:1D-TFLite 4:2D-Keras 5:2D-TFLite 6:Ens-Keras 7:Ens-TFLite
# and rows: 1=model headers, 2=Keras/TFLite headers, 3=Precision, 4=Accuracy, 5=Recall, 6=F1 Score
$tbl = $s.Shapes.Item(25).Table

# Row 3 was mislabeled "Precision" -- it actually holds the Accuracy numbers.
# Retype the label as "Accuracy" and fill in the 1D CNN results.
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Accuracy"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Font.Name = "+mn-lt"
$tbl.Cell(3, 2).Shape.TextFrame.TextRange.Text = "97.9%"
$tbl.Cell(3, 3).Shape.TextFrame.TextRange.Text = "98.8%"

# Row 4 was mislabeled "Accuracy" -- it actually holds the Precision numbers.
$tbl.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Precision"
$tbl.Cell(4, 2).Shape.TextFrame.TextRange.Text = "95.6%"
$tbl.Cell(4, 3).Shape.TextFrame.TextRange.Text = "96.5%"

# Row 5 (Recall) and Row 6 (F1 Score): fill in the 1D CNN results.
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "86.9%"
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "93.8%"
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "91.0%"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "95.1%"
